# "added test data file" - populate Sheet1 with a simple username/password
# test-data table (A1:B2) that a Selenium login test would read from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Data row
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"

# Widen column A a bit so the header/values are readable
$ws.Columns("A").ColumnWidth = 9.166666666666666

# Leave the cursor on B2, matching where the author finished typing
$ws.Range("B2").Select()
